$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing score for 2025-02-04 / activity / Summerbody25 (E12: 0 -> 1) ---
$ws.Range("E12").Value = $true

# --- Adjust column A width (25.6640625 chars -> 25 chars) ---
$ws.Columns.Item(1).ColumnWidth = 24.17

# --- Append the new daily scores for 2025-02-05 (rows 14-16) ---
# Force column A to be read as text first, so "2025-02-05" is stored as a
# shared string rather than being auto-converted into a date serial number.
$ws.Range("A14:A16").NumberFormat = "@"

# Row 14: sleep
$ws.Range("A14").Value = "2025-02-05"
$ws.Range("B14").Value = "sleep"
$ws.Range("C14").Value = $true
$ws.Range("D14").Value = $false
$ws.Range("E14").Value = $true
$ws.Range("F14").Value = $false
$ws.Range("G14").Value = $false
$ws.Range("H14").Value = $false
$ws.Range("I14").Value = $true
$ws.Range("J14").Value = $true
$ws.Range("K14").Value = $false
$ws.Range("L14").Value = $true
$ws.Range("M14").Value = $true
$ws.Range("N14").Value = $true
$ws.Range("O14").Value = $true

# Row 15: activity
$ws.Range("A15").Value = "2025-02-05"
$ws.Range("B15").Value = "activity"
$ws.Range("C15").Value = $true
$ws.Range("D15").Value = $false
$ws.Range("E15").Value = $true
$ws.Range("F15").Value = $true
$ws.Range("G15").Value = $false
$ws.Range("H15").Value = $true
$ws.Range("I15").Value = $false
$ws.Range("J15").Value = $true
$ws.Range("K15").Value = $true
$ws.Range("L15").Value = $true
$ws.Range("M15").Value = $false
$ws.Range("N15").Value = $false
$ws.Range("O15").Value = $false

# Row 16: weekly_activity
$ws.Range("A16").Value = "2025-02-05"
$ws.Range("B16").Value = "weekly_activity"
$ws.Range("C16").Value = $true
$ws.Range("D16").Value = $false
$ws.Range("E16").Value = $true
$ws.Range("F16").Value = $true
$ws.Range("G16").Value = $false
$ws.Range("H16").Value = $false
$ws.Range("I16").Value = $true
$ws.Range("J16").Value = $true
$ws.Range("K16").Value = $false
$ws.Range("L16").Value = $true
$ws.Range("M16").Value = $true
$ws.Range("N16").Value = $false
$ws.Range("O16").Value = $false

# Remove the auto-applied "text" number format style so the new date cells
# end up with no explicit style, matching plain data rows above them.
$ws.Range("A14:O16").ClearFormats()

# Update the selection to mirror the author's next empty block, and let the
# worksheet dimension recompute from the populated range.
[void]$ws.Range("A17:O19").Select()
